$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 2 (existing data rows 2-20 shift down to 3-21).
$ws.Rows.Item(2).Insert()

# The freshly inserted row has no explicit cell style; copy the border/format
# from the row immediately below (row 3, which holds the old row-2 data and
# still carries the original "bordered" style) onto the new row 2.
$ws.Range("A3:F3").Copy()
$ws.Range("A2:F2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Write the new record's values. Values that look like plain numbers
# ("10008651", "1") must still be stored as *text* (shared-string) cells,
# matching the rest of the sheet (every other numeric-looking cell, e.g.
# column D/E's "1".."6", is stored as t="s"). Assigning a literal value
# directly gets auto-coerced to a numeric cell, so instead we briefly hold
# each value as a quoted-string formula (always text) and then flatten it
# back down to a plain literal via copy / paste-values, which preserves the
# cell's existing style (no quote-prefix, no new number format).
$ws.Range("A2").Formula = "=""10008651"""
$ws.Range("B2").Formula = "=""S/G CHK.NGT ORGNL400"""
$ws.Range("C2").Formula = "=""RCS03N"""
$ws.Range("D2").Formula = "=""1"""
$ws.Range("E2").Formula = "=""1"""
$ws.Range("F2").Formula = "=""RT,(E-1B)"""

$ws.Range("A2:F2").Copy()
$ws.Range("A2:F2").PasteSpecial(-4163)
$excel.CutCopyMode = $false
